$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '60.409.92'
$ws.Range("E2").Value = '  +5.81%  '
$ws.Range("D3").Value = '3.274.01'
$ws.Range("E3").Value = '  +0.89%  '
$ws.Range("E4").Value = '  -0.02%  '
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '406.20'
$ws.Range("E5").Value = '  +2.87%  '
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '110.81'
$ws.Range("E6").Value = '  +2.82%  '
$ws.Range("D7").Value = '3.268.64'
$ws.Range("E7").Value = '  +0.95%  '
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = '0.563'
$ws.Range("E8").Value = '  -4.01%  '
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = '1.00'
$ws.Range("E9").Value = '  +0.00%  '
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '0.615'
$ws.Range("E10").Value = '  -1.28%  '
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '0.113'
$ws.Range("E11").Value = '  +13.83%  '
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = '38.34'
$ws.Range("E12").Value = '  -1.86%  '
$ws.Range("E13").Value = '  -0.22%  '
$ws.Range("D14").Value = '3.797.95'
$ws.Range("E14").Value = '  +0.98%  '
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = '8.08'
$ws.Range("E15").Value = '  -1.79%  '
$ws.Range("E16").Value = '  -1.52%  '
$ws.Range("D17").Value = '3.279.85'
$ws.Range("E17").Value = '  +1.23%  '
$ws.Range("D18").Value = '60.379.77'
$ws.Range("E18").Value = '  +6.01%  '
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '0.977'
$ws.Range("E19").Value = '  -5.22%  '
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '10.47'
$ws.Range("E20").Value = '  -2.77%  '
$ws.Range("E21").Value = '  +1.73%  '
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '3.25'
$ws.Range("E22").Value = '  -2.13%  '
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '12.39'
$ws.Range("E23").Value = '  -3.91%  '
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '294.99'
$ws.Range("E24").Value = '  -0.21%  '
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '72.86'
$ws.Range("E25").Value = '  -1.78%  '
$ws.Range("E26").Value = '  -3.08%  '
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '28.83'
$ws.Range("E27").Value = '  +3.12%  '
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = '4.26'
$ws.Range("E28").Value = '  -2.33%  '
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = '0.171'
$ws.Range("E29").Value = '  +2.33%  '
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '7.32'
$ws.Range("E30").Value = '  +1.23%  '
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '7.43'
$ws.Range("E31").Value = '  -1.93%  '
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '0.112'
$ws.Range("E32").Value = '  +2.87%  '
$ws.Range("E33").Value = '  +0.10%  '
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '11.07'
$ws.Range("E34").Value = '  -2.04%  '
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = '2.44'
$ws.Range("E35").Value = '  +14.60%  '
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = '38.45'
$ws.Range("E36").Value = '  -1.88%  '
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = '0.0476'
$ws.Range("E37").Value = '  -0.92%  '
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '52.10'
$ws.Range("E38").Value = '  +1.20%  '
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '0.999'
$ws.Range("E39").Value = '  -0.04%  '
$ws.Range("E40").Value = '  +5.73%  '
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '3.27'
$ws.Range("E41").Value = '  -5.88%  '
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '134.46'
$ws.Range("E42").Value = '  -0.15%  '
$ws.Range("B43").Value = 'TheGraph'
$ws.Range("C43").Value = 'https://coinranking.com/coin/qhd1biQ7M+thegraph-grt'
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '0.287'
$ws.Range("E43").Value = '  +2.72%  '
$ws.Range("B44").Value = 'Stellar'
$ws.Range("C44").Value = 'https://coinranking.com/coin/f3iaFeCKEmkaZ+stellar-xlm'
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '0.119'
$ws.Range("E44").Value = '  -3.06%  '
$ws.Range("E45").Value = '  -1.12%  '
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '16.11'
$ws.Range("E46").Value = '  -5.22%  '
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '3.72'
$ws.Range("E47").Value = '  -5.38%  '
$ws.Range("E48").Value = '  +2.33%  '
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '20.75'
$ws.Range("E49").Value = '  -5.97%  '
$ws.Range("D50").Value = '2.101.06'
$ws.Range("E50").Value = '  -2.46%  '
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '2.35'
$ws.Range("E51").Value = '  +0.31%  '
